$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 384.09000000000003
$ws.Range("C2").Value = 0.9742857142857142
$ws.Range("D2").Value = 0.9399999999999993
$ws.Range("E2").Value = 0.9880000000000001
$ws.Range("F2").Value = 0.9542087542087544
$ws.Range("G2").Value = 0.9979999999999999
$ws.Range("H2").Value = 0.9399999999999993
$ws.Range("I2").Value = 0.9797979797979798
$ws.Range("J2").Value = 0.9305626598465471
$ws.Range("K2").Value = 205.77118286331955
$ws.Range("L2").Value = 0.05877888856275234
$ws.Range("M2").Value = 0.1780931546050317
$ws.Range("N2").Value = 0.04773665131519141
$ws.Range("O2").Value = 0.1780931546050317
$ws.Range("P2").Value = 0.07994019615852753
$ws.Range("Q2").Value = 0.10802135702432011
$ws.Range("R2").Value = 0.16724484273754414
$ws.Range("S2").Value = 0.014070529413622775
$ws.Range("B3").Value = 348.5400000000001
$ws.Range("C3").Value = 0.9728571428571425
$ws.Range("D3").Value = 0.9399999999999995
$ws.Range("E3").Value = 0.9860000000000001
$ws.Range("F3").Value = 0.9460000000000001
$ws.Range("G3").Value = 0.9979999999999999
$ws.Range("H3").Value = 0.9399999999999995
$ws.Range("I3").Value = 0.9766666666666667
$ws.Range("J3").Value = 0.9292838874680304
$ws.Range("K3").Value = 191.4755790531455
$ws.Range("L3").Value = 0.056325320629098084
$ws.Range("M3").Value = 0.16329931618554522
$ws.Range("N3").Value = 0.05128647999525317
$ws.Range("O3").Value = 0.16329931618554522
$ws.Range("P3").Value = 0.0854774666587465
$ws.Range("Q3").Value = 0.11554584263763874
$ws.Range("R3").Value = 0.14849134973506964
$ws.Range("S3").Value = 0.014070529413622775
$ws.Range("B4").Value = 308.0600000000001
$ws.Range("C4").Value = 0.9742857142857142
$ws.Range("D4").Value = 0.9399999999999997
$ws.Range("E4").Value = 0.9880000000000001
$ws.Range("F4").Value = 0.954208754208754
$ws.Range("G4").Value = 0.9969999999999997
$ws.Range("H4").Value = 0.9399999999999997
$ws.Range("I4").Value = 0.9797979797979798
$ws.Range("J4").Value = 0.9305626598465473
$ws.Range("K4").Value = 175.2514972053799
$ws.Range("L4").Value = 0.05877888856275234
$ws.Range("M4").Value = 0.1780931546050317
$ws.Range("N4").Value = 0.04773665131519141
$ws.Range("O4").Value = 0.1780931546050317
$ws.Range("P4").Value = 0.07994019615852753
$ws.Range("Q4").Value = 0.10802135702431945
$ws.Range("R4").Value = 0.167244842737545
$ws.Range("S4").Value = 0.017144660799776435
$ws.Range("B5").Value = 271.5700000000001
$ws.Range("C5").Value = 0.9714285714285711
$ws.Range("D5").Value = 0.9399999999999997
$ws.Range("E5").Value = 0.9839999999999999
$ws.Range("F5").Value = 0.9505050505050504
$ws.Range("G5").Value = 0.9969999999999997
$ws.Range("H5").Value = 0.9399999999999997
$ws.Range("I5").Value = 0.9713804713804717
$ws.Range("J5").Value = 0.924636828644501
$ws.Range("K5").Value = 162.56965904006822
$ws.Range("L5").Value = 0.06420951071078772
$ws.Range("M5").Value = 0.1780931546050317
$ws.Range("N5").Value = 0.054531984868860406
$ws.Range("O5").Value = 0.1780931546050317
$ws.Range("P5").Value = 0.09828577208032914
$ws.Range("Q5").Value = 0.11461892262280252
$ws.Range("R5").Value = 0.17685193236449603
$ws.Range("S5").Value = 0.017144660799776435
$ws.Range("B6").Value = 237.39
$ws.Range("C6").Value = 0.9742857142857142
$ws.Range("D6").Value = 0.9399999999999995
$ws.Range("E6").Value = 0.9880000000000001
$ws.Range("F6").Value = 0.9542087542087542
$ws.Range("G6").Value = 0.9969999999999997
$ws.Range("H6").Value = 0.9399999999999995
$ws.Range("I6").Value = 0.9797979797979798
$ws.Range("J6").Value = 0.9305626598465468
$ws.Range("K6").Value = 149.4908459749393
$ws.Range("L6").Value = 0.05877888856275112
$ws.Range("M6").Value = 0.1780931546050317
$ws.Range("N6").Value = 0.04773665131518991
$ws.Range("O6").Value = 0.1780931546050317
$ws.Range("P6").Value = 0.07994019615852481
$ws.Range("Q6").Value = 0.10802135702432011
$ws.Range("R6").Value = 0.16724484273754414
$ws.Range("S6").Value = 0.017144660799776435
$ws.Range("B7").Value = 201.11
$ws.Range("C7").Value = 0.9742857142857142
$ws.Range("D7").Value = 0.9399999999999995
$ws.Range("E7").Value = 0.9880000000000001
$ws.Range("F7").Value = 0.9542087542087542
$ws.Range("H7").Value = 0.9399999999999995
$ws.Range("I7").Value = 0.9797979797979798
$ws.Range("J7").Value = 0.930562659846547
$ws.Range("K7").Value = 135.77111570060055
$ws.Range("L7").Value = 0.05877888856275112
$ws.Range("M7").Value = 0.1780931546050317
$ws.Range("O7").Value = 0.1780931546050317
$ws.Range("P7").Value = 0.07994019615852481
$ws.Range("Q7").Value = 0.10802135702432011
$ws.Range("R7").Value = 0.16724484273754414
$ws.Range("S7").Value = 0.017144660799776435
$ws.Range("B8").Value = 158.12000000000006
$ws.Range("C8").Value = 0.9714285714285713
$ws.Range("D8").Value = 0.9349999999999997
$ws.Range("E8").Value = 0.9860000000000001
$ws.Range("F8").Value = 0.9491582491582491
$ws.Range("H8").Value = 0.9349999999999997
$ws.Range("I8").Value = 0.9747474747474748
$ws.Range("J8").Value = 0.923562659846547
$ws.Range("K8").Value = 106.32210257020643
$ws.Range("L8").Value = 0.06420951071078436
$ws.Range("M8").Value = 0.18333333333333335
$ws.Range("O8").Value = 0.18333333333333335
$ws.Range("P8").Value = 0.09332493998835673
$ws.Range("Q8").Value = 0.11716060979917028
$ws.Range("R8").Value = 0.17857462092446869
$ws.Range("S8").Value = 0.017144660799776435
$ws.Range("B9").Value = 108.69000000000001
$ws.Range("C9").Value = 0.9728571428571428
$ws.Range("D9").Value = 0.9399999999999995
$ws.Range("F9").Value = 0.946333333333333
$ws.Range("H9").Value = 0.9399999999999995
$ws.Range("I9").Value = 0.9749999999999996
$ws.Range("J9").Value = 0.9294450127877237
$ws.Range("K9").Value = 87.59776010696606
$ws.Range("L9").Value = 0.059873432092561664
$ws.Range("M9").Value = 0.16329931618554522
$ws.Range("O9").Value = 0.16329931618554522
$ws.Range("Q9").Value = 0.11994153083608858
$ws.Range("R9").Value = 0.1561214522212931
$ws.Range("S9").Value = 0.017144660799776435
$ws.Range("B10").Value = 69.44000000000004
$ws.Range("C10").Value = 0.9714285714285713
$ws.Range("D10").Value = 0.9399999999999997
$ws.Range("E10").Value = 0.9839999999999999
$ws.Range("F10").Value = 0.9439999999999998
$ws.Range("G10").Value = 0.9965
$ws.Range("H10").Value = 0.9399999999999997
$ws.Range("I10").Value = 0.9733333333333333
$ws.Range("J10").Value = 0.9262404092071607
$ws.Range("K10").Value = 66.37806014555464
$ws.Range("L10").Value = 0.05743073230052987
$ws.Range("M10").Value = 0.16329931618554522
$ws.Range("N10").Value = 0.054531984868860406
$ws.Range("O10").Value = 0.16329931618554522
$ws.Range("P10").Value = 0.09088664144809612
$ws.Range("Q10").Value = 0.11632996145731278
$ws.Range("R10").Value = 0.15013714762777752
$ws.Range("S10").Value = 0.017773831632741082
$ws.Range("B11").Value = 25.240000000000006
$ws.Range("C11").Value = 0.9728571428571431
$ws.Range("D11").Value = 0.9349999999999996
$ws.Range("E11").Value = 0.9880000000000001
$ws.Range("F11").Value = 0.9446666666666665
$ws.Range("G11").Value = 0.9904999999999998
$ws.Range("H11").Value = 0.9349999999999996
$ws.Range("I11").Value = 0.98
$ws.Range("J11").Value = 0.9282097186700768
$ws.Range("K11").Value = 35.61100726438159
$ws.Range("L11").Value = 0.05632532062909936
$ws.Range("M11").Value = 0.16899883449481556
$ws.Range("N11").Value = 0.04773665131518991
$ws.Range("O11").Value = 0.16899883449481556
$ws.Range("P11").Value = 0.079561085525313
$ws.Range("Q11").Value = 0.11799104000679657
$ws.Range("R11").Value = 0.15057242276567973
$ws.Range("S11").Value = 0.04643611882266878
